# Fill the missing Q4-2015 quarter for OMI CLEAR: insert a new data row right
# after the header row, push the existing quarters down by one row, and
# populate the new row with the Q4-2015 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a brand-new row above the current row 2 (the first data row). This
# shifts all existing data rows (2-11) down to (3-12) and keeps their values
# and number formats intact.
$ws.Rows.Item(2).Insert()

# Populate the newly-inserted row 2 with the Q4-2015 figures.
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "OMI CLEAR"
$ws.Cells.Item(2, 3).Value = "Q4-2015"
$ws.Cells.Item(2, 4).Value = "aggregated"
$ws.Cells.Item(2, 5).Value = 1941895.15
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 27081905.04
$ws.Cells.Item(2, 9).Value = 27081905.04
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 8525718
$ws.Cells.Item(2, 12).Value = 27081905.04
$ws.Cells.Item(2, 13).Value = "For additional information regarding participants commitment to replenish the default fund, please see OMIClear Instruction B07-2014."
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 39333212.3
$ws.Cells.Item(2, 16).Value = 26468203.3
$ws.Cells.Item(2, 17).Value = 40506167.31
$ws.Cells.Item(2, 18).Value = 63336290.38
$ws.Cells.Item(2, 19).Value = 191774701.12
$ws.Cells.Item(2, 20).Value = 191774701.12
$ws.Cells.Item(2, 21).Value = 6
$ws.Cells.Item(2, 22).Value = 12
$ws.Cells.Item(2, 23).Value = "n.a."
$ws.Cells.Item(2, 24).Value = 67
$ws.Cells.Item(2, 25).Value = 5
$ws.Cells.Item(2, 26).Value = 0
$ws.Cells.Item(2, 27).Value = 25353291.04
$ws.Cells.Item(2, 28).Value = 0
$ws.Cells.Item(2, 29).Value = 700518106.15

# Match the style used by the other rows' "CCP index" column (A).
$ws.Cells.Item(2, 1).Style = $ws.Cells.Item(3, 1).Style

# Re-number the sequential index column (A) for every data row (0,1,2,...)
# now that the new row has been inserted.
$lastRow = $ws.Cells.Item(1, 1).CurrentRegion.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
